# daily auto push: 2026-02-25 03:15 UTC
#
# A new reading was captured for 2026/02/25 (水) at time-rank 8, value 48.
# It belongs right after the existing 2026/02/25 row (row 853), so insert a
# new row 854 and push every following row down by one. The sheet's used
# range grows from A1:D895 to A1:D896 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 854 (and everything below it) down by one row.
$ws.Rows.Item(854).Insert()

# Row 853 ("2026/02/25", "水", 5, 48) is the closest match to the new entry
# ("2026/02/25", "水", 8, 48) - copy it into the freshly inserted row so the
# date/weekday text cells keep their original text representation (typing
# "2026/02/25" directly would be auto-recognised as a date value), then fix
# up the one column that actually differs.
$ws.Range("A853:D853").Copy($ws.Range("A854:D854"))
$ws.Range("C854").Value = 8

Write-Output "inserted row 854: 2026/02/25 / 水 / 8 / 48"
